$d = $word.ActiveDocument

# Locate the target paragraph: the "From Bill:  Here's an idea: ..." item.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "From Bill:*") {
        $target = $cand
        break
    }
}

$pRange = $target.Range
$paraStart = $pRange.Start
$paraEnd = $pRange.End

# Range covering just the run content, excluding the trailing paragraph mark.
$contentRange = $d.Range($paraStart, $paraEnd - 1)
$contentRange.Text = ""

$newRunsXml = '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr><w:t>I did a first attempt at the following' + [char]0x2026 + '</w:t></w:r>' + `
'<w:proofErr w:type="gramStart"/>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr><w:t>.</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>F</w:t></w:r>' + `
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>rom</w:t></w:r>' + `
'<w:proofErr w:type="gramEnd"/>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve"> Bill:  </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve">Here' + [char]0x2019 + 's an idea: how about when doing a budget entry for a subaccount, you show the description in addition to the abbreviation.' + [char]0x00A0 + ' E.g. PHY201 ' + [char]0x2013 + ' Astronomy.' + [char]0x00A0 + ' That' + [char]0x2019 + 's for newbies like me who don' + [char]0x2019 + 't know the courses by heart J.' + [char]0x00A0 + ' Maybe also a </w:t></w:r>' + `
'<w:proofErr w:type="spellStart"/>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>mouseover</w:t></w:r>' + `
'<w:proofErr w:type="spellEnd"/>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve"> on the summary pages?</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'

$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($paraStart, $paraStart)
$insertionPoint.InsertXML($pkg)
